# Generate Report for Handback
# The "da6fb246-739e-42e4-b222-7ad388577717" file has now been handed back
# (in sync with en-US). This moves its row ahead of the still-pending
# "54780db5-304c-4ec0-95e5-c0ffd5f74adc" row on every sheet, flips its
# status to "Handed back: in sync with en-US", and fills in the
# Latest Target File / Latest Handback File / Latest Handback DateTime
# columns that only apply once a handback has happened.

$wb = $excel.ActiveWorkbook

$daFile  = "da6fb246-739e-42e4-b222-7ad388577717.md"
$wuFile  = "54780db5-304c-4ec0-95e5-c0ffd5f74adc.md"

$daUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/7628b386195e18d119695ed96dc6ae4dae487866/e2e/da6fb246-739e-42e4-b222-7ad388577717.md"
$wuUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/77df2c771df368776fc5dcf21734e6843cb02678/e2e/54780db5-304c-4ec0-95e5-c0ffd5f74adc.md"

$daZhXlf = "da6fb246-739e-42e4-b222-7ad388577717.3d0f7469a98155ceac1e3a7bcd5d25c99444cd79.zh-cn.xlf"
$wuZhXlf = "54780db5-304c-4ec0-95e5-c0ffd5f74adc.f92723e7ff041355eed3a5ee80a5b97a6e51a368.zh-cn.xlf"
$daZhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8d03c18f70dc4c3b9d24b27a5a539bd0d737922/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/da6fb246-739e-42e4-b222-7ad388577717.3d0f7469a98155ceac1e3a7bcd5d25c99444cd79.zh-cn.xlf"
$wuZhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbc4abfea99c085cb49035fa45ff075028c894b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/54780db5-304c-4ec0-95e5-c0ffd5f74adc.f92723e7ff041355eed3a5ee80a5b97a6e51a368.zh-cn.xlf"

$daDeXlf = "da6fb246-739e-42e4-b222-7ad388577717.3d0f7469a98155ceac1e3a7bcd5d25c99444cd79.de-de.xlf"
$wuDeXlf = "54780db5-304c-4ec0-95e5-c0ffd5f74adc.f92723e7ff041355eed3a5ee80a5b97a6e51a368.de-de.xlf"
$daDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/153ed50e28fab17efacdbeaebc0d9d2a9587922a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/da6fb246-739e-42e4-b222-7ad388577717.3d0f7469a98155ceac1e3a7bcd5d25c99444cd79.de-de.xlf"
$wuDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/88357aa6711acee491ede9504217029c19b9c60a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/54780db5-304c-4ec0-95e5-c0ffd5f74adc.f92723e7ff041355eed3a5ee80a5b97a6e51a368.de-de.xlf"

$handedBack = "Handed back: in sync with en-US"
$readyForHandoff = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet "Overview": da6fb246 row moves to row 2, 54780db5 row to row 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A1").Hyperlinks.Delete()

$wsOverview.Range("A2").Value = $daFile
$wsOverview.Range("B2").Value = $handedBack
$wsOverview.Range("C2").Value = $handedBack
$wsOverview.Range("D2").Value = "2016-28-20 02:28:28"

$wsOverview.Range("A3").Value = $wuFile
$wsOverview.Range("B3").Value = $readyForHandoff
$wsOverview.Range("C3").Value = $readyForHandoff
$wsOverview.Range("D3").Value = "2016-28-20 02:28:02"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $daUrl, "", "", $daFile) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $wuUrl, "", "", $wuFile) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": same row swap, plus populate Latest Target File /
# Latest Handback File / Latest Handback DateTime for the handed-back row
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A1").Hyperlinks.Delete()

$wsZh.Range("A2").Value = $daFile
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $handedBack
$wsZh.Range("D2").Value = $daZhXlf
$wsZh.Range("E2").Value = "2016-03-20 02:28:25"
$wsZh.Range("F2").Value = $daFile
$wsZh.Range("G2").Value = $daZhXlf
$wsZh.Range("H2").Value = "2016-03-20 02:28:45"
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = $wuFile
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $readyForHandoff
$wsZh.Range("D3").Value = $wuZhXlf
$wsZh.Range("E3").Value = "2016-03-20 02:27:59"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $daUrl, "", "", $daFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $daUrl, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $daZhXlfUrl, "", "", $daZhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $daUrl, "", "", $daFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $daZhXlfUrl, "", "", $daZhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $wuUrl, "", "", $wuFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $wuUrl, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $wuZhXlfUrl, "", "", $wuZhXlf) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": same row swap, plus populate Latest Target File /
# Latest Handback File / Latest Handback DateTime for the handed-back row
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A1").Hyperlinks.Delete()

$wsDe.Range("A2").Value = $daFile
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $handedBack
$wsDe.Range("D2").Value = $daDeXlf
$wsDe.Range("E2").Value = "2016-03-20 02:28:28"
$wsDe.Range("F2").Value = $daFile
$wsDe.Range("G2").Value = $daDeXlf
$wsDe.Range("H2").Value = "2016-03-20 02:28:50"
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = $wuFile
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $readyForHandoff
$wsDe.Range("D3").Value = $wuDeXlf
$wsDe.Range("E3").Value = "2016-03-20 02:28:02"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $daUrl, "", "", $daFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $daUrl, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $daDeXlfUrl, "", "", $daDeXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $daUrl, "", "", $daFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $daDeXlfUrl, "", "", $daDeXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $wuUrl, "", "", $wuFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $wuUrl, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $wuDeXlfUrl, "", "", $wuDeXlf) | Out-Null
